$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Item" names (column A) to more descriptive variants
$ws.Range("A2").Value = "Agulha de talagarça"
$ws.Range("A3").Value = "Agenda escolar"
$ws.Range("A4").Value = "Apontador escolar"
$ws.Range("A6").Value = "Bloco Criativo papel colorido A4"
$ws.Range("A7").Value = "Borracha escolar"
$ws.Range("A8").Value = "Caderno de desenho"
$ws.Range("A10").Value = "Cartolina folha"
$ws.Range("A11").Value = "Cola Bastão escolar"
$ws.Range("A12").Value = "Cola Branca escolar"
$ws.Range("A13").Value = "Estojo escolar"
$ws.Range("A15").Value = "Lápis de Cor caixa 12 ou 24 cores"
$ws.Range("A16").Value = "Lápis Preto n 2"
$ws.Range("A17").Value = "Massa de Modelar escolar"
$ws.Range("A18").Value = "Papel A4 pacote"
$ws.Range("A19").Value = "Pasta escolar"
$ws.Range("A20").Value = "Régua escolar"
$ws.Range("A21").Value = "Tesoura escolar"
$ws.Range("A22").Value = "Tinta Guache pote 250 gr"
$ws.Range("A23").Value = "TNT metro"
$ws.Range("A24").Value = "Pasta escolar"
$ws.Range("A25").Value = "Pincel n 10 "
$ws.Range("A26").Value = "Régua escolar"
$ws.Range("A28").Value = "Tesoura escolar"
$ws.Range("A29").Value = "Tinta Guache pote 250 gr"
$ws.Range("A30").Value = "TNT metro"

# Update the sheet view: scroll position and active cell/selection
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A30").Select()
